$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 348.0909
$ws.Range("I6").Value = 193.25
$ws.Range("K6").Value = 579.75
$ws.Range("M6").Value = -467.75
$ws.Range("H32").Value = 20702.5
$ws.Range("J32").Value = 21043.6
$ws.Range("L32").Value = 21043.6
$ws.Range("N32").Value = -21695.6
$ws.Range("H40").Value = 6391.76
$ws.Range("J40").Value = 9642.643
$ws.Range("L40").Value = 9642.643
$ws.Range("N40").Value = -9992.643
$ws.Range("H111").Value = 1421.2222
$ws.Range("I111").Value = 1473.875
$ws.Range("K111").Value = 4421.625
$ws.Range("M111").Value = -1354.625
$ws.Range("H137").Value = 2777.4644
$ws.Range("I137").Value = 1969.2
$ws.Range("J137").Value = 4798.125
$ws.Range("K137").Value = 5907.6
$ws.Range("L137").Value = 14394.375
$ws.Range("M137").Value = -3357.6
$ws.Range("N137").Value = -19494.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 76927624
$ws.Range("I45").Value = 142859150
$ws.Range("K45").Value = 142859150
$ws.Range("M45").Value = -142858773
$ws.Range("H97").Value = 1622.7142
$ws.Range("I97").Value = 1598.625
$ws.Range("K97").Value = 1598.625
$ws.Range("M97").Value = -1102.625
$ws.Range("H102").Value = 1881.5
$ws.Range("I102").Value = 1881.5
$ws.Range("K102").Value = 1881.5
$ws.Range("M102").Value = -259.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value = 48500
$ws.Range("J63").Value = 48500
$ws.Range("L63").Value = 48500
$ws.Range("N63").Value = -49872
$ws.Range("H66").Value = 48500
$ws.Range("J66").Value = 48500
$ws.Range("L66").Value = 145500
$ws.Range("N66").Value = -152364
$ws.Range("H68").Value = 40000
$ws.Range("J68").Value = 40000
$ws.Range("L68").Value = 40000
$ws.Range("H71").Value = 40000
$ws.Range("J71").Value = 40000
$ws.Range("L71").Value = 120000
$ws.Range("H105").Value = 5511.56
$ws.Range("I105").Value = 1232.2727
$ws.Range("J105").Value = 8873.857
$ws.Range("K105").Value = 1232.2727
$ws.Range("L105").Value = 8873.857
$ws.Range("M105").Value = 514.7273
$ws.Range("N105").Value = -12367.857
$ws.Range("N68").Value = -41622
$ws.Range("N71").Value = -128112

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 22405.557
$ws.Range("I31").Value = 3726.1843
$ws.Range("K31").Value = 3726.1843
$ws.Range("M31").Value = -3431.1843
$ws.Range("H34").Value = 22405.557
$ws.Range("I34").Value = 3726.1843
$ws.Range("K34").Value = 3726.1843
$ws.Range("M34").Value = -3524.1843
$ws.Range("H74").Value = 600000
$ws.Range("J74").Value = 600000
$ws.Range("L74").Value = 600000
$ws.Range("H77").Value = 600000
$ws.Range("J77").Value = 600000
$ws.Range("L77").Value = 1800000
$ws.Range("H81").Value = 39000
$ws.Range("J81").Value = 39000
$ws.Range("L81").Value = 39000
$ws.Range("H82").Value = 38996
$ws.Range("J82").Value = 38996
$ws.Range("L82").Value = 38996
$ws.Range("N82").Value = -39718
$ws.Range("H84").Value = 39000
$ws.Range("J84").Value = 39000
$ws.Range("L84").Value = 117000
$ws.Range("H85").Value = 38996
$ws.Range("J85").Value = 38996
$ws.Range("L85").Value = 38996
$ws.Range("N85").Value = -41492
$ws.Range("H134").Value = 2697
$ws.Range("I134").Value = 1652.7142
$ws.Range("J134").Value = 10007
$ws.Range("K134").Value = 4958.142599999999
$ws.Range("L134").Value = 30021
$ws.Range("M134").Value = -2423.142599999999
$ws.Range("N134").Value = -35091
$ws.Range("N74").Value = -601748
$ws.Range("N77").Value = -1808736
$ws.Range("N81").Value = -40996
$ws.Range("N84").Value = -126984

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 11176.643
$ws.Range("I5").Value = 491.375
$ws.Range("J5").Value = 25423.666
$ws.Range("K5").Value = 1474.125
$ws.Range("L5").Value = 76270.99800000001
$ws.Range("M5").Value = -1362.125
$ws.Range("N5").Value = -76494.99800000001
$ws.Range("H7").Value = 38851.92
$ws.Range("I7").Value = 732
$ws.Range("J7").Value = 62676.875
$ws.Range("K7").Value = 2196
$ws.Range("L7").Value = 188030.625
$ws.Range("M7").Value = -2084
$ws.Range("N7").Value = -188254.625
$ws.Range("H22").Value = 3543
$ws.Range("I22").Value = 1390
$ws.Range("K22").Value = 4170
$ws.Range("M22").Value = -4001
$ws.Range("H27").Value = 3543
$ws.Range("I27").Value = 1390
$ws.Range("K27").Value = 4170
$ws.Range("M27").Value = -4068
$ws.Range("H32").Value = 10499.7
$ws.Range("J32").Value = 10499.7
$ws.Range("L32").Value = 31499.1
$ws.Range("N32").Value = -32065.1
$ws.Range("H39").Value = 4000
$ws.Range("I39").Value = 4000
$ws.Range("J39").Value = 4000
$ws.Range("K39").Value = 12000
$ws.Range("L39").Value = 12000
$ws.Range("M39").Value = -11706
$ws.Range("N39").Value = -12588
$ws.Range("H113").Value = 1091.0834
$ws.Range("I113").Value = 1264.6666
$ws.Range("K113").Value = 3793.9998
$ws.Range("M113").Value = -1623.9998
$ws.Range("H135").Value = 11176.643
$ws.Range("I135").Value = 491.375
$ws.Range("J135").Value = 25423.666
$ws.Range("K135").Value = 4422.375
$ws.Range("L135").Value = 228812.994
$ws.Range("M135").Value = -1887.375
$ws.Range("N135").Value = -233882.994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5505.364
$ws.Range("I126").Value = 1773.75
$ws.Range("J126").Value = 7637.7144
$ws.Range("K126").Value = 5321.25
$ws.Range("L126").Value = 22913.1432
$ws.Range("M126").Value = -2851.25
$ws.Range("N126").Value = -27853.1432
$ws.Range("H132").Value = 46129.117
$ws.Range("I132").Value = 54349.668
$ws.Range("J132").Value = 11602.8
$ws.Range("K132").Value = 163049.004
$ws.Range("L132").Value = 34808.39999999999
$ws.Range("M132").Value = -160519.004
$ws.Range("N132").Value = -39868.39999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 5000
$ws.Range("I42").Value = 5000
$ws.Range("K42").Value = 5000
$ws.Range("M42").Value = -4437
$ws.Range("H49").Value = 5000
$ws.Range("I49").Value = 5000
$ws.Range("K49").Value = 5000
$ws.Range("M49").Value = -4853
$ws.Range("H61").Value = 9342.190000000001
$ws.Range("I61").Value = 7749.5
$ws.Range("J61").Value = 12527.571
$ws.Range("K61").Value = 7749.5
$ws.Range("L61").Value = 12527.571
$ws.Range("M61").Value = -7547.5
$ws.Range("N61").Value = -12931.571
$ws.Range("H62").Value = 490000
$ws.Range("J62").Value = 490000
$ws.Range("L62").Value = 490000
$ws.Range("N62").Value = -491248
$ws.Range("H64").Value = 26333.334
$ws.Range("J64").Value = 34500
$ws.Range("L64").Value = 34500
$ws.Range("N64").Value = -34950
$ws.Range("H65").Value = 490000
$ws.Range("J65").Value = 490000
$ws.Range("L65").Value = 1470000
$ws.Range("N65").Value = -1476240
$ws.Range("H67").Value = 26333.334
$ws.Range("J67").Value = 34500
$ws.Range("L67").Value = 34500
$ws.Range("N67").Value = -36060
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H93").Value = 3180.4375
$ws.Range("I93").Value = 2768.7693
$ws.Range("K93").Value = 2768.7693
$ws.Range("M93").Value = -1520.7693
$ws.Range("H113").Value = 9342.190000000001
$ws.Range("I113").Value = 7749.5
$ws.Range("J113").Value = 12527.571
$ws.Range("K113").Value = 7749.5
$ws.Range("L113").Value = 12527.571
$ws.Range("M113").Value = -5579.5
$ws.Range("N113").Value = -16867.571
$ws.Range("H136").Value = 3795.5833
$ws.Range("I136").Value = 1449.2222
$ws.Range("J136").Value = 10834.667
$ws.Range("K136").Value = 4347.6666
$ws.Range("L136").Value = 32504.001
$ws.Range("M136").Value = -1797.6666
$ws.Range("N136").Value = -37604.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4472.136
$ws.Range("I136").Value = 3286.1187
$ws.Range("J136").Value = 14468.571
$ws.Range("K136").Value = 9858.356100000001
$ws.Range("L136").Value = 43405.713
$ws.Range("M136").Value = -7308.356100000001
$ws.Range("N136").Value = -48505.713
